$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ZonaComun")

# Insert new row at 8 (shifts old row 8 "idTipoZonaComun" down to row 9)
$ws.Rows.Item(8).Insert()

# Copy formatting from row 9 (the shifted original row) into new row 8
$ws.Range("A9:U9").Copy()
$ws.Range("A8:U8").PasteSpecial(-4122)

# Set new row 8 content: new attribute "conjuntoResidencial" of type "ConjuntoResidencial"
$ws.Range("A8").Value = "conjuntoResidencial"
$ws.Range("B8").Value = "ConjuntoResidencial"

# Rename old attribute from "idTipoZonaComun" to "tipoZonaComun" (now on row 9)
$ws.Range("A9").Value = "tipoZonaComun"

# Fix up hyperlinks: delete all and re-add with corrected addresses (rows >= 9 shift down by 1)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A1"), "", "'Objetos de Dominio'!A1", "", "Volver al inicio")
$ws.Hyperlinks.Add($ws.Range("I20"), "", "'Tipo Relación Institución'!A6", "", "'Tipo Relación Institución'!A6")
$ws.Hyperlinks.Add($ws.Range("S4"), "", "'Objeto Dominio 2'!A17", "", "'Objeto Dominio 2'!A17")
$ws.Hyperlinks.Add($ws.Range("T4"), "", "'Objeto Dominio 2'!A18", "", "'Objeto Dominio 2'!A18")
$ws.Hyperlinks.Add($ws.Range("U4"), "", "'Objeto Dominio 2'!A19", "", "'Objeto Dominio 2'!A19")
$ws.Hyperlinks.Add($ws.Range("A18:B18"), "", "'Objeto Dominio 2'!R4", "", "Reponsabilidad 2")
$ws.Hyperlinks.Add($ws.Range("A17:B17"), "", "'Objeto Dominio 2'!Q4", "", "Reponsabilidad 1")
$ws.Hyperlinks.Add($ws.Range("A20:B20"), "", "'Objeto Dominio 2'!T4", "", "Reponsabilidad 4")
$ws.Hyperlinks.Add($ws.Range("R4"), "", "'Objeto Dominio 2'!A16", "", "'Objeto Dominio 2'!A16")
$ws.Hyperlinks.Add($ws.Range("A1:Q1"), "", "'Listado Objetos de Dominio'!A1", "", "<-Volver al inicio")
$ws.Hyperlinks.Add($ws.Range("A19:B19"), "", "'Objeto Dominio 2'!S4", "", "Reponsabilidad 3")
$ws.Hyperlinks.Add($ws.Range("C13"), "", "ZonaComun!A8", "", "identificador")

Write-Host "done"
